$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (COM ColumnWidth adds ~0.8333 padding vs the stored
# XML "width" attribute, so subtract that offset to land on the exact
# target stored widths of 24 / 23 / 25).
$ws.Columns(1).ColumnWidth = 23.1666666666667
$ws.Columns(2).ColumnWidth = 22.1666666666667
$ws.Columns(3).ColumnWidth = 24.1666666666667

# Update header row values
$ws.Range("A1").Value = "var_1_input_first_name"
$ws.Range("B1").Value = "var_2_input_last_name"
$ws.Range("C1").Value = "var_3_input_postal_code"

# Update data row values
$ws.Range("A2").Value = "amandarobinson"
$ws.Range("B2").Value = "sextondeanna"
$ws.Range("C2").Value = "pamelapatterson"
